$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: TRON -> Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.73"
$ws.Range("E13").Value = "  -3.45%  "

# Row 14: Polkadot -> TRON
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.100"
$ws.Range("E14").Value = "  -3.27%  "

$ws.Range("D2").Value = "41.695.41"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.172.52"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.31"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.62"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.09"
$ws.Range("E10").Value = "  -6.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -5.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.39"
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D15").Value = "2.497.84"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.32"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "2.169.18"
$ws.Range("E17").Value = "  -4.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.779"
$ws.Range("E18").Value = "  -7.22%  "
$ws.Range("D19").Value = "41.572.09"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.04"
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  -6.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.96"
$ws.Range("E23").Value = "  -12.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.38"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -5.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.26"
$ws.Range("E28").Value = "  -9.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.17"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.02"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.84"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.40"
$ws.Range("E33").Value = "  +7.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0774"
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.121"
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.31"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -7.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0310"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.15"
$ws.Range("E40").Value = "  -8.38%  "
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.38"
$ws.Range("E42").Value = "  -5.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.93"
$ws.Range("E43").Value = "  -9.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.44"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.189"
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0967"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.64"
$ws.Range("E47").Value = "  -6.82%  "
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.12"
$ws.Range("E49").Value = "  -4.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  -7.05%  "
$ws.Range("E51").Value = "  -2.28%  "
